$d = $word.ActiveDocument

# 1) <meta name="title" content="...">
#    "1 Introduction to Game Maker" -> "9 Exiting the Room"
$r = $d.Content
$found = $r.Find.Execute("1 Introduction to Game Maker")
if ($found) {
    $r.Text = "9 Exiting the Room"
}

# 2) <meta name="description" content="...">
#    Remove the yellow highlight and replace the whole sentence with the
#    new tutorial description.
$r = $d.Content
$found = $r.Find.Execute("This is our first article in a new series of Game Maker, where we will be introducing you to it. ")
if ($found) {
    $r.Text = "In this tutorial, we will be learning how we can get the hero to exit the room. "
    $r.HighlightColorIndex = 0
}

# 3) <meta name="revised" content="...">
#    "Thursday 30th, 2025" (with superscript "th") -> " November 29, 2025"
$r = $d.Content
$found = $r.Find.Execute("Thursday 30th, 2025")
if ($found) {
    $r.Text = " November 29, 2025"
}

# 4) <meta name="url" content="...">
#    "I-Snuck-A-Book/PDF_Optimizer.html" -> new article path (with trailing space)
$r = $d.Content
$found = $r.Find.Execute("I-Snuck-A-Book/PDF_Optimizer.html")
if ($found) {
    $r.Text = "Enlightenment/Articles/2025/4_Game_Maker/9_Exiting_the_Room/9_Exiting_the_Room.html "
}
